# =====================================================================
# Rebuild workbook: rename Sheet1, insert Order Week column, add 3 sheets
# =====================================================================
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet1 -> "Sales vs PO": insert "Order Week" column (new col C)
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Sales vs PO"

# Insert a new column at C; the existing PO_Requested_Qty column moves to D
$ws1.Columns.Item(3).Insert()

# Give the new column C the same date number-format/style as column A
$ws1.Range("A2:A29").Copy()
$ws1.Range("C2:C29").PasteSpecial(-4122)
$ws1.Application.CutCopyMode = $false

$ws1.Range("C1").Value = "Order Week"

# Column C (Order Week) receives the original "ds" dates;
# column A (ds) is shifted forward; column D (PO_Requested_Qty) becomes 0
$ws1.Cells.Item(2, 3).Value = 45460
$ws1.Cells.Item(2, 1).Value = 45466
$ws1.Cells.Item(2, 4).Value = 0
$ws1.Cells.Item(3, 3).Value = 45467
$ws1.Cells.Item(3, 1).Value = 45473
$ws1.Cells.Item(3, 4).Value = 0
$ws1.Cells.Item(4, 3).Value = 45474
$ws1.Cells.Item(4, 1).Value = 45480
$ws1.Cells.Item(4, 4).Value = 0
$ws1.Cells.Item(5, 3).Value = 45481
$ws1.Cells.Item(5, 1).Value = 45487
$ws1.Cells.Item(5, 4).Value = 0
$ws1.Cells.Item(6, 3).Value = 45488
$ws1.Cells.Item(6, 1).Value = 45494
$ws1.Cells.Item(6, 4).Value = 0
$ws1.Cells.Item(7, 3).Value = 45495
$ws1.Cells.Item(7, 1).Value = 45501
$ws1.Cells.Item(7, 4).Value = 0
$ws1.Cells.Item(8, 3).Value = 45502
$ws1.Cells.Item(8, 1).Value = 45508
$ws1.Cells.Item(8, 4).Value = 0
$ws1.Cells.Item(9, 3).Value = 45509
$ws1.Cells.Item(9, 1).Value = 45515
$ws1.Cells.Item(9, 4).Value = 0
$ws1.Cells.Item(10, 3).Value = 45516
$ws1.Cells.Item(10, 1).Value = 45522
$ws1.Cells.Item(10, 4).Value = 0
$ws1.Cells.Item(11, 3).Value = 45523
$ws1.Cells.Item(11, 1).Value = 45529
$ws1.Cells.Item(11, 4).Value = 0
$ws1.Cells.Item(12, 3).Value = 45530
$ws1.Cells.Item(12, 1).Value = 45536
$ws1.Cells.Item(12, 4).Value = 0
$ws1.Cells.Item(13, 3).Value = 45537
$ws1.Cells.Item(13, 1).Value = 45543
$ws1.Cells.Item(13, 4).Value = 0
$ws1.Cells.Item(14, 3).Value = 45544
$ws1.Cells.Item(14, 1).Value = 45550
$ws1.Cells.Item(14, 4).Value = 0
$ws1.Cells.Item(15, 3).Value = 45551
$ws1.Cells.Item(15, 1).Value = 45557
$ws1.Cells.Item(15, 4).Value = 0
$ws1.Cells.Item(16, 3).Value = 45558
$ws1.Cells.Item(16, 1).Value = 45564
$ws1.Cells.Item(16, 4).Value = 0
$ws1.Cells.Item(17, 3).Value = 45565
$ws1.Cells.Item(17, 1).Value = 45571
$ws1.Cells.Item(17, 4).Value = 0
$ws1.Cells.Item(18, 3).Value = 45572
$ws1.Cells.Item(18, 1).Value = 45578
$ws1.Cells.Item(18, 4).Value = 0
$ws1.Cells.Item(19, 3).Value = 45579
$ws1.Cells.Item(19, 1).Value = 45585
$ws1.Cells.Item(19, 4).Value = 0
$ws1.Cells.Item(20, 3).Value = 45586
$ws1.Cells.Item(20, 1).Value = 45592
$ws1.Cells.Item(20, 4).Value = 0
$ws1.Cells.Item(21, 3).Value = 45593
$ws1.Cells.Item(21, 1).Value = 45599
$ws1.Cells.Item(21, 4).Value = 0
$ws1.Cells.Item(22, 3).Value = 45600
$ws1.Cells.Item(22, 1).Value = 45606
$ws1.Cells.Item(22, 4).Value = 0
$ws1.Cells.Item(23, 3).Value = 45607
$ws1.Cells.Item(23, 1).Value = 45613
$ws1.Cells.Item(23, 4).Value = 0
$ws1.Cells.Item(24, 3).Value = 45614
$ws1.Cells.Item(24, 1).Value = 45620
$ws1.Cells.Item(24, 4).Value = 0
$ws1.Cells.Item(25, 3).Value = 45621
$ws1.Cells.Item(25, 1).Value = 45627
$ws1.Cells.Item(25, 4).Value = 0
$ws1.Cells.Item(26, 3).Value = 45628
$ws1.Cells.Item(26, 1).Value = 45634
$ws1.Cells.Item(26, 4).Value = 0
$ws1.Cells.Item(27, 3).Value = 45635
$ws1.Cells.Item(27, 1).Value = 45641
$ws1.Cells.Item(27, 4).Value = 0
$ws1.Cells.Item(28, 3).Value = 45642
$ws1.Cells.Item(28, 1).Value = 45648
$ws1.Cells.Item(28, 4).Value = 0
$ws1.Cells.Item(29, 3).Value = 45649
$ws1.Cells.Item(29, 1).Value = 45655
$ws1.Cells.Item(29, 4).Value = 0

# ---------------------------------------------------------------
# Add "Weekly Growth", "Volume Insights", "Prediction Info" sheets
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Weekly Growth"

$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Volume Insights"

$ws4 = $wb.Worksheets.Add($null, $ws3)
$ws4.Name = "Prediction Info"

# ---------------------------------------------------------------
# Apply matching header / date styles (copied from Sheet1) before
# writing values, so the same style indexes are reused
# ---------------------------------------------------------------
$ws1.Range("A1:C1").Copy()
$ws2.Range("A1:C1").PasteSpecial(-4122)
$ws1.Application.CutCopyMode = $false

$ws1.Range("A2:A15").Copy()
$ws2.Range("A2:A15").PasteSpecial(-4122)
$ws1.Application.CutCopyMode = $false

$ws1.Range("A1:C1").Copy()
$ws3.Range("A1:D1").PasteSpecial(-4122)
$ws1.Application.CutCopyMode = $false

$ws1.Range("A1").Copy()
$ws4.Range("A1").PasteSpecial(-4122)
$ws1.Application.CutCopyMode = $false

# ---------------------------------------------------------------
# Weekly Growth data
# ---------------------------------------------------------------
$ws2.Range("A1").Value = "ds"
$ws2.Range("B1").Value = "PO_Requested_Qty"
$ws2.Range("C1").Value = "Growth%"
$ws2.Cells.Item(2, 1).Value = 45467
$ws2.Cells.Item(2, 2).Value = 20
$ws2.Cells.Item(2, 3).Value = 0
$ws2.Cells.Item(3, 1).Value = 45481
$ws2.Cells.Item(3, 2).Value = 80
$ws2.Cells.Item(3, 3).Value = 300
$ws2.Cells.Item(4, 1).Value = 45488
$ws2.Cells.Item(4, 2).Value = 80
$ws2.Cells.Item(4, 3).Value = 0
$ws2.Cells.Item(5, 1).Value = 45502
$ws2.Cells.Item(5, 2).Value = 120
$ws2.Cells.Item(5, 3).Value = 50
$ws2.Cells.Item(6, 1).Value = 45509
$ws2.Cells.Item(6, 2).Value = 280
$ws2.Cells.Item(6, 3).Value = 133.3333333333333
$ws2.Cells.Item(7, 1).Value = 45516
$ws2.Cells.Item(7, 2).Value = 120
$ws2.Cells.Item(7, 3).Value = -57.14285714285714
$ws2.Cells.Item(8, 1).Value = 45530
$ws2.Cells.Item(8, 2).Value = 200
$ws2.Cells.Item(8, 3).Value = 66.66666666666667
$ws2.Cells.Item(9, 1).Value = 45537
$ws2.Cells.Item(9, 2).Value = 580
$ws2.Cells.Item(9, 3).Value = 190
$ws2.Cells.Item(10, 1).Value = 45544
$ws2.Cells.Item(10, 2).Value = 220
$ws2.Cells.Item(10, 3).Value = -62.06896551724138
$ws2.Cells.Item(11, 1).Value = 45551
$ws2.Cells.Item(11, 2).Value = 1340
$ws2.Cells.Item(11, 3).Value = 509.0909090909091
$ws2.Cells.Item(12, 1).Value = 45579
$ws2.Cells.Item(12, 2).Value = 20
$ws2.Cells.Item(12, 3).Value = -98.50746268656717
$ws2.Cells.Item(13, 1).Value = 45586
$ws2.Cells.Item(13, 2).Value = 700
$ws2.Cells.Item(13, 3).Value = 3400
$ws2.Cells.Item(14, 1).Value = 45593
$ws2.Cells.Item(14, 2).Value = 320
$ws2.Cells.Item(14, 3).Value = -54.28571428571429
$ws2.Cells.Item(15, 1).Value = 45607
$ws2.Cells.Item(15, 2).Value = 40
$ws2.Cells.Item(15, 3).Value = -87.5

# ---------------------------------------------------------------
# Volume Insights data
# ---------------------------------------------------------------
$ws3.Range("A1").Value = "Total_PO_Quantity"
$ws3.Range("B1").Value = "Average_PO_Quantity"
$ws3.Range("C1").Value = "Max_PO_Quantity"
$ws3.Range("D1").Value = "Min_PO_Quantity"
$ws3.Cells.Item(2, 1).Value = 4120
$ws3.Cells.Item(2, 2).Value = 294.2857142857143
$ws3.Cells.Item(2, 3).Value = 1340
$ws3.Cells.Item(2, 4).Value = 20

# ---------------------------------------------------------------
# Prediction Info data
# ---------------------------------------------------------------
$ws4.Range("A1").Value = "Predicted_Next_Week_PO_Quantity"
$ws4.Cells.Item(2, 1).Value = 521.098901098901

# Return to the first sheet as the active sheet/selection
$ws1.Activate()
$ws1.Range("A1").Select()
